# Updated symbol list on Mon Jan  2 13:42:42 UTC 2023 with GitHub Actions
# Refresh of coin price/volume figures (and a 3-row reshuffle of the
# KickToken / BKEXToken / CEJI block at rows 41-43).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value, whether the value is numeric-looking
# text (Price / Volume columns) that must be written with a leading quote
# so Excel stores it as literal text (matching the workbook's existing
# inline-string cells) instead of auto-converting it to a Number/Percent.
$updates = @(
    @{ Cell = 'D2'; Value = '246.61'; AsText = $true }
    @{ Cell = 'E2'; Value = '0.95%'; AsText = $true }
    @{ Cell = 'D3'; Value = '29.85'; AsText = $true }
    @{ Cell = 'E3'; Value = '9.54%'; AsText = $true }
    @{ Cell = 'D4'; Value = '5.163'; AsText = $true }
    @{ Cell = 'E4'; Value = '0.94%'; AsText = $true }
    @{ Cell = 'D5'; Value = '0.05712'; AsText = $true }
    @{ Cell = 'E5'; Value = '0.96%'; AsText = $true }
    @{ Cell = 'D6'; Value = '6.608'; AsText = $true }
    @{ Cell = 'E6'; Value = '2.05%'; AsText = $true }
    @{ Cell = 'D7'; Value = '0.8586'; AsText = $true }
    @{ Cell = 'E7'; Value = '4.52%'; AsText = $true }
    @{ Cell = 'D8'; Value = '0.8672'; AsText = $true }
    @{ Cell = 'E8'; Value = '3.20%'; AsText = $true }
    @{ Cell = 'E9'; Value = '2.67%'; AsText = $true }
    @{ Cell = 'D10'; Value = '0.07081'; AsText = $true }
    @{ Cell = 'E10'; Value = '2.20%'; AsText = $true }
    @{ Cell = 'D11'; Value = '0.02861'; AsText = $true }
    @{ Cell = 'E11'; Value = '-4.27%'; AsText = $true }
    @{ Cell = 'E12'; Value = '0.04%'; AsText = $true }
    @{ Cell = 'D13'; Value = '0.001521'; AsText = $true }
    @{ Cell = 'E13'; Value = '-0.24%'; AsText = $true }
    @{ Cell = 'D14'; Value = '0.04144'; AsText = $true }
    @{ Cell = 'E14'; Value = '-1.49%'; AsText = $true }
    @{ Cell = 'D15'; Value = '0.0005992'; AsText = $true }
    @{ Cell = 'E15'; Value = '0.17%'; AsText = $true }
    @{ Cell = 'D16'; Value = '0.006194'; AsText = $true }
    @{ Cell = 'E16'; Value = '0.91%'; AsText = $true }
    @{ Cell = 'E17'; Value = '3,764.33%'; AsText = $true }
    @{ Cell = 'E18'; Value = '-0.99%'; AsText = $true }
    @{ Cell = 'D19'; Value = '3.057'; AsText = $true }
    @{ Cell = 'E19'; Value = '1.75%'; AsText = $true }
    @{ Cell = 'D20'; Value = '2.173'; AsText = $true }
    @{ Cell = 'E20'; Value = '-5.85%'; AsText = $true }
    @{ Cell = 'D21'; Value = '0.3146'; AsText = $true }
    @{ Cell = 'E21'; Value = '1.05%'; AsText = $true }
    @{ Cell = 'D22'; Value = '0.03250'; AsText = $true }
    @{ Cell = 'E22'; Value = '3.52%'; AsText = $true }
    @{ Cell = 'E23'; Value = '0.71%'; AsText = $true }
    @{ Cell = 'D24'; Value = '3.508'; AsText = $true }
    @{ Cell = 'E24'; Value = '-1.51%'; AsText = $true }
    @{ Cell = 'E25'; Value = '0.45%'; AsText = $true }
    @{ Cell = 'D26'; Value = '0.005096'; AsText = $true }
    @{ Cell = 'E26'; Value = '14.26%'; AsText = $true }
    @{ Cell = 'E27'; Value = '-0.21%'; AsText = $true }
    @{ Cell = 'D28'; Value = '0.0001211'; AsText = $true }
    @{ Cell = 'E28'; Value = '23.50%'; AsText = $true }
    @{ Cell = 'D40'; Value = '0.03771'; AsText = $true }
    @{ Cell = 'E40'; Value = '3.25%'; AsText = $true }
    @{ Cell = 'B41'; Value = 'BKEXToken'; AsText = $false }
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'; AsText = $false }
    @{ Cell = 'D41'; Value = '0.1071'; AsText = $true }
    @{ Cell = 'E41'; Value = '1.77%'; AsText = $true }
    @{ Cell = 'B42'; Value = 'CEJI'; AsText = $false }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'; AsText = $false }
    @{ Cell = 'D42'; Value = '0.002602'; AsText = $true }
    @{ Cell = 'E42'; Value = '13.08%'; AsText = $true }
    @{ Cell = 'B43'; Value = 'KickToken'; AsText = $false }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'; AsText = $false }
    @{ Cell = 'D43'; Value = '0.005674'; AsText = $true }
    @{ Cell = 'E43'; Value = '-5.95%'; AsText = $true }
    @{ Cell = 'D44'; Value = '0.009349'; AsText = $true }
    @{ Cell = 'E44'; Value = '4.14%'; AsText = $true }
    @{ Cell = 'D45'; Value = '0.00005092'; AsText = $true }
    @{ Cell = 'E45'; Value = '-4.03%'; AsText = $true }
    @{ Cell = 'E46'; Value = '-0.01%'; AsText = $true }
    @{ Cell = 'D47'; Value = '0.07512'; AsText = $true }
    @{ Cell = 'E47'; Value = '-25.64%'; AsText = $true }
    @{ Cell = 'D48'; Value = '0.002733'; AsText = $true }
    @{ Cell = 'E48'; Value = '7.04%'; AsText = $true }
    @{ Cell = 'D49'; Value = '0.00002101'; AsText = $true }
    @{ Cell = 'E49'; Value = '-0.01%'; AsText = $true }
    @{ Cell = 'E50'; Value = '-0.01%'; AsText = $true }
)

foreach ($u in $updates) {
    $value = $u.Value
    if ($u.AsText) {
        # Leading apostrophe = Excel's quote-prefix: forces text storage
        # while the apostrophe itself is stripped from the stored value.
        $value = "'" + $value
    }
    $ws.Range($u.Cell).Value = $value
}

